$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.648.06'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').Value = '2.810.42'
$ws.Range('E3').Value = '  +1.42%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '352.18'
$ws.Range('E5').Value = '  +5.48%  '
$ws.Range('D6').Value = '112.89'
$ws.Range('E6').Value = '  -3.40%  '
$ws.Range('E7').Value = '  +3.82%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.595'
$ws.Range('E9').Value = '  +3.21%  '
$ws.Range('D10').Value = '41.26'
$ws.Range('E10').Value = '  -1.92%  '
$ws.Range('E11').Value = '  -1.72%  '
$ws.Range('D13').Value = '19.86'
$ws.Range('E13').Value = '  -2.24%  '
$ws.Range('D14').Value = '7.70'
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').Value = '3.251.16'
$ws.Range('E15').Value = '  +1.40%  '
$ws.Range('D16').Value = '2.813.98'
$ws.Range('E16').Value = '  +1.01%  '
$ws.Range('D17').Value = '0.884'
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('D18').Value = '51.360.20'
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('D19').Value = '7.43'
$ws.Range('E19').Value = '  +8.10%  '
$ws.Range('D20').Value = '3.18'
$ws.Range('E20').Value = '  -4.28%  '
$ws.Range('D21').Value = '13.30'
$ws.Range('E21').Value = '  -1.69%  '
$ws.Range('E22').Value = '  +1.21%  '
$ws.Range('D23').Value = '270.25'
$ws.Range('E23').Value = '  -3.04%  '
$ws.Range('D24').Value = '69.48'
$ws.Range('E24').Value = '  -0.57%  '
$ws.Range('D25').Value = '2.74'
$ws.Range('E25').Value = '  +1.84%  '
$ws.Range('D26').Value = '26.62'
$ws.Range('E26').Value = '  -0.90%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Value = '10.27'
$ws.Range('E28').Value = '  +0.80%  '
$ws.Range('E29').Value = '  +0.45%  '
$ws.Range('E30').Value = '  -2.41%  '
$ws.Range('D31').Value = '50.55'
$ws.Range('E31').Value = '  +0.58%  '
$ws.Range('D32').Value = '33.83'
$ws.Range('E32').Value = '  -3.93%  '
$ws.Range('D33').Value = '5.83'
$ws.Range('E33').Value = '  +4.39%  '
$ws.Range('D34').Value = '0.0445'
$ws.Range('E34').Value = '  +25.34%  '
$ws.Range('D35').Value = '0.0820'
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '5.08'
$ws.Range('E36').Value = '  +0.78%  '
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('E38').Value = '  -2.03%  '
$ws.Range('D39').Value = '3.20'
$ws.Range('E39').Value = '  -1.42%  '
$ws.Range('D40').Value = '18.04'
$ws.Range('E40').Value = '  -6.40%  '
$ws.Range('D41').Value = '23.53'
$ws.Range('E41').Value = '  +0.40%  '
$ws.Range('E42').Value = '  +1.65%  '
$ws.Range('D43').Value = '126.41'
$ws.Range('E43').Value = '  -1.62%  '
$ws.Range('E44').Value = '  +1.45%  '
$ws.Range('E45').Value = '  -1.18%  '
$ws.Range('D46').Value = '2.071.25'
$ws.Range('E46').Value = '  -0.97%  '
$ws.Range('D47').Value = '3.30'
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('D48').Value = '2.29'
$ws.Range('E48').Value = '  +2.46%  '
$ws.Range('D49').Value = '5.63'
$ws.Range('E49').Value = '  +1.41%  '
$ws.Range('D50').Value = '0.921'
$ws.Range('E50').Value = '  +4.78%  '
$ws.Range('D51').Value = '60.54'
$ws.Range('E51').Value = '  +0.01%  '
